# Auto-generated edit script applying the Goblin_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice* / LeveProfit* columns (H-N)
# for specific leve rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 983.125
$ws.Range("I17").Value = 440
$ws.Range("J17").Value = 1019.3333
$ws.Range("K17").Value = 1320
$ws.Range("L17").Value = 3057.9999
$ws.Range("M17").Value = -1152
$ws.Range("N17").Value = -3393.9999
$ws.Range("H88").Value = 4521.385
$ws.Range("J88").Value = 5689.8
$ws.Range("L88").Value = 5689.8
$ws.Range("N88").Value = -6501.8
$ws.Range("H91").Value = 4521.385
$ws.Range("J91").Value = 5689.8
$ws.Range("L91").Value = 5689.8
$ws.Range("N91").Value = -8497.799999999999
$ws.Range("H113").Value = 4181.5
$ws.Range("I113").Value = 3561
$ws.Range("J113").Value = 4802
$ws.Range("K113").Value = 3561
$ws.Range("L113").Value = 4802
$ws.Range("M113").Value = -307
$ws.Range("N113").Value = -11310
$ws.Range("H116").Value = 4984.1665
$ws.Range("I116").Value = 4702.5
$ws.Range("K116").Value = 4702.5
$ws.Range("M116").Value = -1260.5
$ws.Range("H132").Value = 2223.2942
$ws.Range("I132").Value = 1655.3182
$ws.Range("J132").Value = 3264.5833
$ws.Range("K132").Value = 4965.9546
$ws.Range("L132").Value = 9793.749899999999
$ws.Range("M132").Value = -2435.9546
$ws.Range("N132").Value = -14853.7499
$ws.Range("H137").Value = 3942.7715
$ws.Range("I137").Value = 4406.5864
$ws.Range("K137").Value = 13219.7592
$ws.Range("M137").Value = -10669.7592

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2586
$ws.Range("I61").Value = 2675
$ws.Range("J61").Value = 1829.5
$ws.Range("K61").Value = 2675
$ws.Range("L61").Value = 1829.5
$ws.Range("M61").Value = -2463
$ws.Range("N61").Value = -2253.5
$ws.Range("H97").Value = 1548.2307
$ws.Range("I97").Value = 1702.7
$ws.Range("K97").Value = 1702.7
$ws.Range("M97").Value = -1206.7
$ws.Range("H122").Value = 4446847
$ws.Range("I122").Value = 6175120
$ws.Range("J122").Value = 2716.4285
$ws.Range("K122").Value = 18525360
$ws.Range("L122").Value = 8149.2855
$ws.Range("M122").Value = -18522910
$ws.Range("N122").Value = -13049.2855
$ws.Range("H136").Value = 2586
$ws.Range("I136").Value = 2675
$ws.Range("J136").Value = 1829.5
$ws.Range("K136").Value = 8025
$ws.Range("L136").Value = 5488.5
$ws.Range("M136").Value = -5475
$ws.Range("N136").Value = -10588.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 3372.5
$ws.Range("I14").Value = 990
$ws.Range("K14").Value = 990
$ws.Range("M14").Value = -818
$ws.Range("H134").Value = 2466.7646
$ws.Range("I134").Value = 1766.7858
$ws.Range("J134").Value = 5733.3335
$ws.Range("K134").Value = 5300.357400000001
$ws.Range("L134").Value = 17200.0005
$ws.Range("M134").Value = -2765.357400000001
$ws.Range("N134").Value = -22270.0005

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 110.333336
$ws.Range("I7").Value = 40.75
$ws.Range("J7").Value = 249.5
$ws.Range("K7").Value = 40.75
$ws.Range("L7").Value = 249.5
$ws.Range("M7").Value = 72.25
$ws.Range("N7").Value = -475.5
$ws.Range("H41").Value = 37216.223
$ws.Range("I41").Value = 5000
$ws.Range("K41").Value = 5000
$ws.Range("M41").Value = -4572
$ws.Range("H108").Value = 40000
$ws.Range("J108").Value = 40000
$ws.Range("L108").Value = 40000
$ws.Range("N108").Value = -47680

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 20367.084
$ws.Range("I74").Value = 19628.75
$ws.Range("K74").Value = 58886.25
$ws.Range("M74").Value = -57825.25
$ws.Range("H77").Value = 20367.084
$ws.Range("I77").Value = 19628.75
$ws.Range("K77").Value = 176658.75
$ws.Range("M77").Value = -171354.75
$ws.Range("H131").Value = 1962855
$ws.Range("J131").Value = 2471507.8
$ws.Range("L131").Value = 7414523.399999999
$ws.Range("N131").Value = -7424603.399999999

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 550
$ws.Range("I5").Value = 550
$ws.Range("K5").Value = 550
$ws.Range("M5").Value = -438
$ws.Range("H68").Value = 75000
$ws.Range("J68").Value = 75000
$ws.Range("L68").Value = 75000
$ws.Range("N68").Value = -76622
$ws.Range("H71").Value = 75000
$ws.Range("J71").Value = 75000
$ws.Range("L71").Value = 225000
$ws.Range("N71").Value = -233112
$ws.Range("H80").Value = 22301.666
$ws.Range("I80").Value = 50905
$ws.Range("J80").Value = 8000
$ws.Range("K80").Value = 50905
$ws.Range("L80").Value = 8000
$ws.Range("M80").Value = -49907
$ws.Range("N80").Value = -9996
$ws.Range("H83").Value = 22301.666
$ws.Range("I83").Value = 50905
$ws.Range("J83").Value = 8000
$ws.Range("K83").Value = 254525
$ws.Range("L83").Value = 40000
$ws.Range("M83").Value = -249533
$ws.Range("N83").Value = -49984
$ws.Range("H97").Value = 859.0833
$ws.Range("J97").Value = 957.6667
$ws.Range("L97").Value = 957.6667
$ws.Range("N97").Value = -1949.6667

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3070.5264
$ws.Range("I22").Value = 2494.2222
$ws.Range("K22").Value = 2494.2222
$ws.Range("M22").Value = -2199.2222
$ws.Range("N22").Value = -4179.2
$ws.Range("H27").Value = 3070.5264
$ws.Range("I27").Value = 2494.2222
$ws.Range("K27").Value = 2494.2222
$ws.Range("M27").Value = -2387.2222
$ws.Range("N27").Value = -3803.2
$ws.Range("H55").Value = 2058.25
$ws.Range("I55").Value = 742.5
$ws.Range("K55").Value = 742.5
$ws.Range("M55").Value = -569.5
$ws.Range("N55").Value = -3720
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H111").Value = 208000
$ws.Range("J111").Value = 208000
$ws.Range("L111").Value = 208000
$ws.Range("N111").Value = -216180
$ws.Range("H116").Value = 271666.66
$ws.Range("J116").Value = 271666.66
$ws.Range("L116").Value = 271666.66
$ws.Range("N116").Value = -280844.66

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 13250
$ws.Range("I34").Value = 1500
$ws.Range("K34").Value = 1500
$ws.Range("M34").Value = -1297
$ws.Range("H40").Value = 14899
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 14899
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 14899
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -15197
$ws.Range("H56").Value = 5142.5
$ws.Range("I56").Value = 5285
$ws.Range("J56").Value = 5000
$ws.Range("K56").Value = 5285
$ws.Range("L56").Value = 5000
$ws.Range("M56").Value = -4571
$ws.Range("N56").Value = -6428
$ws.Range("H132").Value = 5864.8667
$ws.Range("I132").Value = 4596.077
$ws.Range("J132").Value = 6835.1177
$ws.Range("K132").Value = 13788.231
$ws.Range("L132").Value = 20505.3531
$ws.Range("M132").Value = -11258.231
$ws.Range("N132").Value = -25565.3531

Write-Host "Applied Goblin_Profits.xlsx updates"